$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = -13.38399999999999
$ws.Range("C12").Value = -10.8975
$ws.Range("C18").Value = -12.58089999999999
$ws.Range("C37").Value = -13.30499999999999
$ws.Range("C55").Value = -13.4044
$ws.Range("C68").Value = -11.5898
$ws.Range("C77").Value = -11.9115
$ws.Range("C78").Value = -12.1055
